# 第02组/项目计划表.xlsx — "修改用例图" commit
#
# Fills in the second weekly block (rows 13-18) of Sheet1 with the same
# team-member names used in the first block (rows 3-8) and records each
# member's updated task ("修改用例图[...]" / "编写用例描述[...]"). Also
# nudges the saved window tab-ratio and leaves the selection on B18 (the
# last cell touched), matching where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 13: 练富珊 ---------------------------------------------------
$ws.Range("A13").Value = "练富珊"
$ws.Range("B13").Value = "修改用例图[前端手机用户用例图]"

# --- Row 14: 黄成志 ---------------------------------------------------
$ws.Range("A14").Value = "黄成志"
$ws.Range("B14").Value = "编写用例描述[手机端用户注册]"

# --- Row 15: 黄皓燊 ---------------------------------------------------
$ws.Range("A15").Value = "黄皓燊"
$ws.Range("B15").Value = "修改用例图[前端PC用户用例图]"

# --- Row 16: 郑嘉蔚 (task cell gets an explicit font re-stamp, like the
#             author did when touching up this particular entry) ------
$ws.Range("A16").Value = "郑嘉蔚"
$ws.Range("B16").Value = "编写用例描述[手机端用户登陆]"
$ws.Range("B16").Font.Name = "宋体"

# --- Row 17: 郑瑞贤 ---------------------------------------------------
$ws.Range("A17").Value = "郑瑞贤"
$ws.Range("B17").Value = "修改用例图[系统管理员用例图]"

# --- Row 18: 辛伟达 ---------------------------------------------------
$ws.Range("A18").Value = "辛伟达"
$ws.Range("B18").Value = "编写用例描述[手机端创建群]"

# Slightly adjust the saved tab ratio (cosmetic window state: 580 -> 570).
$excel.ActiveWindow.TabRatio = 0.57

# Leave the selection where the author's last edit landed.
$ws.Range("B18").Select()
